# Add 8 new vocabulary rows (14-21) to the "Contracts" sheet, following the
# same left-to-right / row-based data-entry order used by the original author
# so that new shared-string entries land in the same sequence as the target
# workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Contracts")
$ws.Select()

# ---- Row 14: contract ----
$ws.Range("D14").Value = "/ˈkontrӕkt/"
$ws.Range("E14").Value = "hợp đồng"
$ws.Range("F14").Value = "an official written agreement"
$ws.Range("G14").Value = "All employees have a written contract of employment."

# ---- Row 15: integral ----
$ws.Range("A15").Value = "integral"
$ws.Range("D15").Value = "/ˈɪntɪɡrəl/"
$ws.Range("G15").Value = "Music is an integral part of the school's curriculum."
$ws.Range("E15").Value = "một phần không thể thiếu của cái gì đó"
$ws.Range("F15").Value = "[usually before noun] being an essential part of something"

# ---- Row 16: terms ----
$ws.Range("A16").Value = "terms"
$ws.Range("E16").Value = "Điều khoản điều lệ"
$ws.Range("F16").Value = "the conditions that people offer, demand or accept when they make an agreement, an arrangement or a contract"
$ws.Range("G16").Value = "These are the terms and conditions of your employment."
$ws.Range("D16").Value = "/tɜːmz/"

# ---- Row 17: include ----
$ws.Range("A17").Value = "include"
$ws.Range("E17").Value = "bao gồm"
$ws.Range("D17").Value = "/ɪnˈkluːd/"

# ---- Row 18 / 19 / 14 word cells (author filled these out of strict row order) ----
$ws.Range("A18").Value = "certain"
$ws.Range("A19").Value = "quality"
$ws.Range("A14").Value = "contract"

$ws.Range("D18").Value = "/ˈsɜːtn/"
$ws.Range("G18").Value = "I think it was him, but I can't be certain."
$ws.Range("E18").Value = "chính xác, chắc chắn"
$ws.Range("C18").Value = "sure, definite, guaranteed, bound"

$ws.Range("F17").Value = "if one thing includes another, it has the second thing as one of its parts"
$ws.Range("G17").Value = " include something The tour included a visit to the Science Museum."

$ws.Range("F18").Value = "strongly believing something; having no doubts"

$ws.Range("F19").Value = "the standard of something when it is compared to other things like it; how good or bad something is"

# ---- Row 19 / G: rich text (bold lead-in + regular remainder) ----
$ws.Range("G19").Value = "of… quality The soil here is of poor quality"
$r19a = $ws.Range("G19").Characters(1, 11)
$r19a.Font.Name = "Times New Roman"
$r19a.Font.Size = 13
$r19a.Font.Bold = $true
$r19b = $ws.Range("G19").Characters(12, 34)
$r19b.Font.Name = "Times New Roman"
$r19b.Font.Size = 13
$r19b.Font.Bold = $false

# ---- Row 20: deliver ----
$ws.Range("A20").Value = "deliver"
$ws.Range("F20").Value = " to take goods, letters, etc. to the person or people they have been sent to"

# ---- Row 20 / G: rich text (bold lead-in + regular remainder) ----
$ws.Range("G20").Value = " deliver something to somebody/something Leaflets have been delivered to every household."
$r20a = $ws.Range("G20").Characters(1, 40)
$r20a.Font.Name = "Times New Roman"
$r20a.Font.Size = 13
$r20a.Font.Bold = $true
$r20b = $ws.Range("G20").Characters(41, 49)
$r20b.Font.Name = "Times New Roman"
$r20b.Font.Size = 13
$r20b.Font.Bold = $false

$ws.Range("C19").Value = "excellence"
$ws.Range("C20").Value = "save"
$ws.Range("E20").Value = "Giao, trình bày"
$ws.Range("D20").Value = "/dɪˈlɪvə(r)/"

# ---- Row 21: condition ----
$ws.Range("A21").Value = "condition"
$ws.Range("D21").Value = "/kənˈdɪʃn/"
$ws.Range("C21").Value = "disease"
$ws.Range("E21").Value = "Điều kiện, quy định, tình trạng"
$ws.Range("F21").Value = "the state that something is in"

# ---- Row 21 / G: rich text (bold lead-in + regular remainder) ----
$ws.Range("G21").Value = " the condition of something. The condition of the roads is poor."
$r21a = $ws.Range("G21").Characters(1, 28)
$r21a.Font.Name = "Times New Roman"
$r21a.Font.Size = 13
$r21a.Font.Bold = $true
$r21b = $ws.Range("G21").Characters(29, 36)
$r21b.Font.Name = "Times New Roman"
$r21b.Font.Size = 13
$r21b.Font.Bold = $false

# ---- Word forms (column B) reuse existing shared strings ----
$ws.Range("B14").Value = "n"
$ws.Range("B15").Value = "adj"
$ws.Range("B16").Value = "n"
$ws.Range("B17").Value = "v"
$ws.Range("B18").Value = "adj"
$ws.Range("B19").Value = "n"
$ws.Range("B20").Value = "v"
$ws.Range("B21").Value = "n"

$ws.Range("G21").Select()
